# Updated unit test to accommodate change in long term WMO cost
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_modelLogic.py")
$ws.Activate()

# Rows 8-12 and 14-15 (B13 is intentionally left unchanged) get an extra
# *10 multiplier tacked on to their existing "10*TestInputData!B.." formulas,
# reflecting the change in long-term WMO cost being tested.
$rows = @(8, 9, 10, 11, 12, 14, 15)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Formula = $cell.Formula + "*10"
}

# Match the author's final view state: scrolled down so row 19 is at the
# top, with the grand-total cell B41 selected.
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("B41").Select()
